# "add quick copy template"
# - Tidy up the header labels in the shared-string table (drop the
#   double-leading-space padding that was used as a poor-man's indent,
#   collapse "Priority \n(...)" to a single line, de-space "Jira ID",
#   and tweak the SMART-goal remark wording).
# - Move the active selection on Sheet1 from E2 to D4.
# - Populate the previously-empty Sheet2 with a ready-to-copy header row
#   that mirrors Sheet1's header (same styling) so it can be pasted into
#   a fresh tracker quickly.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# 1. Sheet1 header-row text cleanup (shared-string edits)
# ---------------------------------------------------------------------
$ws1.Range("A1").Value = "ID"
$ws1.Range("C1").Value = "Priority(High:1/Medium:2/Low:3)"
$ws1.Range("D1").Value = "Item"
$ws1.Range("E1").Value = "Owner"
$ws1.Range("F1").Value = "Status"
$ws1.Range("G1").Value = "ETA"
$ws1.Range("H1").Value = "Remark"
$ws1.Range("I1").Value = "JiraID"

# SMART-goal template remark: "need to do" -> "requires"
$ws1.Range("D3").Value = "Due to (Why), (Who) requires (How much) of (What) at (Where) in (Scope) by (When), otherwise (Consequence), here is (How) and (Resources) for reference."

# ---------------------------------------------------------------------
# 2. Build the quick-copy header template on Sheet2
# ---------------------------------------------------------------------
$ws1.Range("A1:I1").Copy($ws2.Range("A1:I1"))

$ws2.Rows.Item(1).RowHeight = 28
$ws2.Columns.Item(3).ColumnWidth = 29.5

# ---------------------------------------------------------------------
# 3. Selection bookkeeping: Sheet2 shows the new header highlighted,
#    Sheet1 stays the active sheet/tab with its cursor moved to D4.
# ---------------------------------------------------------------------
$ws2.Range("A1:I1").Select()
$ws1.Range("D4").Select()

Write-Output "quick copy template added to Sheet2"
